$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 250, 50, "Jean", "Aniversário"),
    @(2, 500, 50, "NAejc", "Casamento"),
    @(3, 750, 50, "Jean", "Casamento"),
    @(4, 1000, 50, "NAejc", "Casamento"),
    @(5, 1250, 50, "Jean", "Aniversário"),
    @(6, 1500, 50, "NAejc", "Casamento"),
    @(7, 1750, 50, "Jean", "Aniversário"),
    @(8, 2000, 50, "NAejc", "Casamento"),
    @(9, 2250, 50, "Jean", "Aniversário"),
    @(10, 2500, 50, "NAejc", "Casamento"),
    @(11, 2750, 50, "Jean", "Aniversário"),
    @(12, 3250, 50, "Jean", "Aniversário")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

$ws.Columns.Item(3).ColumnWidth = 14.14

$ws.Range("F4").Select()
